$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "SU changement rotules triangles": update the ball-joint radius
# (Rayon_rotule (mm), row 10, column B) from 9.5 mm to 8 mm.
$ws.Range("B10").Value = 8

# Reflect the author's new scroll/selection position in the sheet view.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B12").Select()
